$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header text changes from "Draft" to "Drafting of manuscript"
$ws.Range("F1").Value = "Drafting of manuscript"

# Update the active selection from E2 to F2
$ws.Range("F2").Select()
